# The author re-uploaded the keyword/appID sheet with the three "helix"
# rows removed (the helix-jump / com.singleton.helix pair at row 4, its
# duplicate at row 7, and the lone "helix" row at row 8). Deleting the
# whole rows shifts everything below up, which also drops the now-unused
# "helix jump", "com.singleton.helix" and "helix" shared strings and
# renumbers the dimension down to A1:B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows("8:8").Delete() | Out-Null
$ws.Rows("7:7").Delete() | Out-Null
$ws.Rows("4:4").Delete() | Out-Null

# The saved selection moves from A5 to A4 once row 4 disappears.
$ws.Range("A4").Select() | Out-Null
